# Refresh the price-comparison table with the latest scrape results.
# The scraper now supports searching items by input from the command
# line, which reshuffled store ordering within each country block and
# picked up several new listings (table grew from 49 to 60 rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number, Country, Shop, Price(BYN), Local price -- one entry per listing.
$rows = @(
    @(2, "RUS", "E-first.ru", 325.59, "9 535 р."),
    @(3, "RUS", "Официальный магазин LG", 324.36, "9 499 р."),
    @(4, "RUS", "Goods.ru", 305.27, "8 940 р."),
    @(5, "RUS", "Kns.ru", 298.65, "8 746 р."),
    @(6, "RUS", "Svyaznoy.ru", 345.56, "10 120 р."),
    @(7, "RUS", "Pleer.ru", 286.12, "8 379 р."),
    @(8, "RUS", "CompYou.ru", 298.27, "8 735 р."),
    @(9, "RUS", "Eldorado.ru", 324.05, "9 490 р."),
    @(10, "RUS", "Elektro-park.ru", 305.95, "8 960 р."),
    @(11, "RUS", "Citilink.ru", 306.98, "8 990 р."),
    @(12, "RUS", "Fotosklad.ru", 317.22, "9 290 р."),
    @(13, "RUS", "ABC.ru", 316.92, "9 281 р."),
    @(14, "RUS", "123.ru", 312.44, "9 150 р."),
    @(15, "RUS", "М.видео", 324.05, "9 490 р."),
    @(16, "RUS", "Just.ru", 310.39, "9 090 р."),
    @(17, "RUS", "Pcplanet.ru", 334.3, "9 790 р."),
    @(18, "RUS", "Kotofoto.ru", 337.03, "9 870 р."),
    @(19, "RUS", "Topcomputer.ru", 285.81, "8 370 р."),
    @(20, "UA", "Itbox.ua", 270.52, "2 999 грн."),
    @(21, "UA", "LuxPRO.ua", 270.52, "2 999 грн."),
    @(22, "UA", "V10.com.ua", 393.02, "4 357 грн."),
    @(23, "UA", "Foroom.com.ua", 306.97, "3 403 грн."),
    @(24, "UA", "Repka.ua", 267.46, "2 965 грн."),
    @(25, "UA", "A-techno.com.ua", 270.52, "2 999 грн."),
    @(26, "UA", "Homebt.com.ua", 286.58, "3 177 грн."),
    @(27, "UA", "Hbox.com.ua", 291.18, "3 228 грн."),
    @(28, "UA", "Brain.com.ua", 270.52, "2 999 грн."),
    @(29, "UA", "Цифра", 262.5, "2 910 грн."),
    @(30, "UA", "Епіцентр", 270.52, "2 999 грн."),
    @(31, "UA", "Denika.ua", 270.52, "2 999 грн."),
    @(32, "UA", "Stylus.ua", 271.52, "3 010 грн."),
    @(33, "UA", "АЛЛО", 270.52, "2 999 грн."),
    @(34, "UA", "Rozetka.ua", 270.52, "2 999 грн."),
    @(35, "BLR", "SOCKET.BY", 283.23, ""),
    @(36, "BLR", "TTN.by", 302.6, ""),
    @(37, "BLR", "KST.by", 283.23, ""),
    @(38, "BLR", "VIPCOMP.BY", 301.35, ""),
    @(39, "BLR", "Sli.by", 283.23, ""),
    @(40, "BLR", "Техник-Сервис", 302.07, ""),
    @(41, "BLR", "Ньютон", 310.0, ""),
    @(42, "BLR", "Медиа Маркет групп ООО", 308.93, ""),
    @(43, "BLR", "24shop.by", 307.25, ""),
    @(44, "BLR", "5 элемент", 359.0, ""),
    @(45, "BLR", "ViP MARKET", 323.29, ""),
    @(46, "BLR", "ITMarket.by", 286.23, ""),
    @(47, "BLR", "bigi", 307.25, ""),
    @(48, "BLR", "ЭЛЕКТРОСИЛА", 326.0, ""),
    @(49, "BLR", "ЧУП `"СВКомп`"", 320.0, ""),
    @(50, "BLR", "RULEZ.BY", 302.07, ""),
    @(51, "BLR", "izliv.by", 301.35, ""),
    @(52, "BLR", "ВИКО-ТЕХНО", 317.0, ""),
    @(53, "BLR", "LevelUP", 317.0, ""),
    @(54, "BLR", "Itsmart.by", 391.32, ""),
    @(55, "BLR", "AMD.by", 290.33, ""),
    @(56, "BLR", "itx.by", 337.87, ""),
    @(57, "BLR", "BITS.By", 317.0, ""),
    @(58, "BLR", "Мультиком.", 310.08, ""),
    @(59, "BLR", "21vek.by", 329.0, ""),
    @(60, "BLR", "imarket.by", 340.51, "")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

